$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.136116072538749
$ws.Range("C2").Value = 0.136116072538749
$ws.Range("D2").Value = 0.775173933911264
$ws.Range("E2").Value = 0.00528136954727986
$ws.Range("F2").Value = 0.5961

$ws.Range("B3").Value = 25.6367580503969
$ws.Range("C3").Value = 0.175594233221896
$ws.Range("E3").Value = 0.99471863045272

$ws.Range("B4").Value = 25.7728741229356
